$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Valid till" value for the second data row (cell C2) by 2 days
$ws.Range("C2").Value = 45716.666666666664

# Move the active selection to C3, matching the saved cursor position
[void]$ws.Range("C3").Select()
